# Update the "Topology Included by Other Algorithms" worksheet (sheet4):
# add a small "Settings:" block (confidence thresholds used for transcript
# abundance contextualization) below the existing table, mirrored in both
# the A:B and D:E column pairs, and widen columns A and D to fit the new
# labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Topology Included by Other Algorithms")

# --- new "Settings:" block (rows 19-22), duplicated in A:B and D:E -------
$ws.Range("A19").Value = "Settings:"
$ws.Range("D19").Value = "Settings:"
$ws.Range("A19").Font.Bold = $true
$ws.Range("D19").Font.Bold = $true

$ws.Range("B20").Value = "Transcript abundance"
$ws.Range("E20").Value = "Transcript abundance"
$ws.Range("B20").HorizontalAlignment = -4108
$ws.Range("E20").HorizontalAlignment = -4108

$ws.Range("A21").Value = "Low confidence: "
$ws.Range("B21").Value = 1
$ws.Range("D21").Value = "Low confidence: "
$ws.Range("E21").Value = 1

$ws.Range("A22").Value = "High confidence:"
$ws.Range("B22").Value = 100
$ws.Range("D22").Value = "High confidence:"
$ws.Range("E22").Value = 100

# --- widen columns A and D to fit the new labels --------------------------
$ws.Columns.Item(1).ColumnWidth = 14.31
$ws.Columns.Item(4).ColumnWidth = 15.69

# --- move the active selection down to reflect the extra rows -------------
$null = $ws.Activate()
$null = $ws.Range("H24").Select()
